# Path to Graduation 2.xlsx - move a few courses around within the
# 2022-2023 block and remove the (empty) "2024" semester block entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# --- Block 1 (Fall/Spring/Summer 2022), rows 4-10 ---
# Row 5 used to carry "CPSC 6180" / 3 in the Summer column (E5:F5).
# Clear those two cells; the course moves down into row 6/7 instead.
$ws.Range("E5:F5").ClearContents()

# Row 6 gains "CPSC 6180" / 3 in the Fall column (A6:B6); the Spring
# column (C6:D6, "CYBR 3115" / 3) was already there and is unchanged.
$ws.Range("A6").Value = "CPSC 6180"
$ws.Range("B6").Value = 3

# New row 7: "CPSC 6185" / 3 in Fall, "CPSC 6985" / 4 in Spring.
$ws.Range("A7").Value = "CPSC 6185"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 6985"
$ws.Range("D7").Value = 4

# --- Block 2 (Fall/Spring/Summer 2023), rows 13-19 ---
# Row 13 used to carry "CPSC 6985" / 4 in the Summer column (E13:F13).
# That course now lives in row 7 above, so clear it here.
$ws.Range("E13:F13").ClearContents()

# Row 15 used to carry "CPSC 6185" / 3 in the Fall column (A15:B15).
# That course now lives in row 7 above, so clear it here.
$ws.Range("A15:B15").ClearContents()

# --- Remove the (unused) Fall/Spring/Summer 2024 block entirely ---
# Row 21 held the 2024 semester headers, row 29 the matching totals row;
# rows 22-28 (the course rows in between) were already blank.
$ws.Range("A21:F21").ClearContents()
$ws.Range("A29:F29").ClearContents()

$wb.Save()
